$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "English" heading above the brief table (NOT the hyperlinked "English"
#    at the very top of the document). Scope the Find to the specific
#    paragraph so the hyperlink occurrence is left untouched.
# ---------------------------------------------------------------------------
$pEnglish = $d.Paragraphs.Item(3)
$pEnglish.Range.Find.Execute("English", $true, $false, $false, $false, $false, $true, 1, $false, "Inglês", 2)

# ---------------------------------------------------------------------------
# 2. Heading text
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Uh oh! Your documents couldn" + [char]8217 + "t be verified", $true, $false, $false, $false, $false, $true, 1, $false, "Uh oh! Os seus documentos não puderam ser verificados", 2)

# ---------------------------------------------------------------------------
# 3. Greeting
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Hi ", $true, $false, $false, $false, $false, $true, 1, $false, "Olá ", 2)

# ---------------------------------------------------------------------------
# 4. Placeholder
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("[PARTNER NAME]", $true, $false, $false, $false, $false, $true, 1, $false, "[NOME DO PARCEIRO]", 2)

# ---------------------------------------------------------------------------
# 5. Intro paragraph
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("We regret to inform you that your documents have failed our verification process as we found the following issues with them: ", $true, $false, $false, $false, $false, $true, 1, $false, "Lamentamos informar que os seus documentos não passaram no nosso processo de verificação, uma vez que encontrámos os seguintes problemas: ", 2)

# ---------------------------------------------------------------------------
# 6. Bullet item bold run
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("A copy of your vaccination certificate", $true, $false, $false, $false, $false, $true, 1, $false, "Uma cópia do seu certificado de vacinação", 2)

# ---------------------------------------------------------------------------
# 7. Bullet item remainder
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(": Document is unclear", $true, $false, $false, $false, $false, $true, 1, $false, ": O documento não é claro", 2)

# ---------------------------------------------------------------------------
# 8. Second bullet placeholder
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("[Document 2]", $true, $false, $false, $false, $false, $true, 1, $false, "[Documento 2]", 2)

# ---------------------------------------------------------------------------
# 9. Resubmission sentence lead-in
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Please resubmit the documents above by ", $true, $false, $false, $false, $false, $true, 1, $false, "Por favor, reenvie os documentos acima até ", 2)

# ---------------------------------------------------------------------------
# 10. Date placeholder
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("DD Mmm YYYY", $true, $false, $false, $false, $false, $true, 1, $false, "DD Mmm AAAA", 2)

# ---------------------------------------------------------------------------
# 11. Resubmission sentence tail
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" so we can proceed with the necessary arrangements.", $true, $false, $false, $false, $false, $true, 1, $false, " para que possamos proceder às devidas diligências.", 2)

# ---------------------------------------------------------------------------
# 12. Contact sentence lead-in. This run sits immediately after a
#     <w:commentRangeStart/> marker. A normal Find/Replace that rewrites the
#     run starting at its very first character causes the marker to be
#     pushed past the run. To keep the marker exactly where it belongs,
#     leave the first character of the run untouched by the bulk edit, then
#     insert the corrected first character and delete the stale one as two
#     separate, narrowly scoped operations.
# ---------------------------------------------------------------------------
$pContact = $d.Paragraphs.Item(21)
$oldLeadIn = "If you have any questions, please contact us via "
$newLeadIn = "Para mais informações, contacte-nos através do "
$leadInStart = $pContact.Range.Start
$bulk = $d.Range($leadInStart + 1, $leadInStart + $oldLeadIn.Length)
$bulk.Text = $newLeadIn.Substring(1)
$insertFirst = $d.Range($leadInStart + 1, $leadInStart + 1)
$insertFirst.InsertBefore($newLeadIn.Substring(0, 1))
$deleteStale = $d.Range($leadInStart, $leadInStart + 1)
$deleteStale.Delete()

# ---------------------------------------------------------------------------
# 13. " or " between the "live chat" and "WhatsApp" hyperlinks -> " ou ".
#     This run has no explicit run formatting in the source; replacing the
#     whole run in one go makes it inherit the neighbouring hyperlink's
#     colour/underline. Leave the leading space untouched (it keeps the
#     run's existing, unformatted identity) and only rewrite "or " -> "ou ".
# ---------------------------------------------------------------------------
$orFind = $pContact.Range.Duplicate()
$orFind.Find.Execute(" or ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$orRange = $d.Range($orFind.Start + 1, $orFind.End)
$orRange.Text = "ou "

# ---------------------------------------------------------------------------
# 14. Country manager sentence lead-in
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("If you have any questions, please contact your country manager, ", $true, $false, $false, $false, $false, $true, 1, $false, "Para mais questões, pode também contactar o seus gestor de parcerias ", 2)

# ---------------------------------------------------------------------------
# 15. ", at " -> ", em "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(", at ", $true, $false, $false, $false, $false, $true, 1, $false, ", em ", 2)

# ---------------------------------------------------------------------------
# 16. Remaining " or " -> " ou " (between [EMAIL ADDRESS] and [WHATSAPP NO]).
#     This one is not adjacent to a hyperlink, so a plain replace is safe.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" or ", $true, $false, $false, $false, $false, $true, 1, $false, " ou ", 2)
